# Daily attendance processing - 2025-11-02 07:20:17
# Reorders the "Recorded By" (column G) values on the Session Analysis
# Results sheet: each comma-separated list of recorders is reversed in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        if ($trimmed.Count -gt 1) {
            $reversed = $trimmed[($trimmed.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversed)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
